$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header cell (G1) onto the new
# "Save" header cell (H1) so it reuses the same bold/bordered/centered style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("H1").Value = "Save"

# Fill the new "Save" column values: first data row is 1, rest are 0.
$ws.Range("H2").Value = 1
$ws.Range("H3:H13").Value = 0
